$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.357.48"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.510.40"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("E11").Value = "  +3.71%  "
$ws.Range("D12").Value = "4.112.20"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "3.513.31"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "64.346.54"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "393.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.576"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("D23").Value = "3.652.49"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("E27").Value = "  +2.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("E32").Value = "  -6.09%  "
$ws.Range("E33").Value = "  +6.61%  "
$ws.Range("D34").Value = "3.541.35"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0788"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.811"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("E46").Value = "  -2.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "2.391.10"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.898"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("E51").Value = "  +0.43%  "
